$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.019430475293696
$ws.Range("D2").Value = 1.024624822111795
$ws.Range("E2").Value = 1.020565154313159
$ws.Range("F2").Value = 1.030511454988552
$ws.Range("I2").Value = 1.028794099063778
$ws.Range("J2").Value = 1.024633284677864
$ws.Range("K2").Value = 1.027452775983998
$ws.Range("L2").Value = 1.023405064274411
$ws.Range("M2").Value = 1.033322254128945

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.020374550968546
$ws.Range("D3").Value = 1.025290576997821
$ws.Range("E3").Value = 1.021365943612955
$ws.Range("F3").Value = 1.031670877349935
$ws.Range("I3").Value = 1.028939656907195
$ws.Range("J3").Value = 1.025214140346149
$ws.Range("K3").Value = 1.02792604620203
$ws.Range("L3").Value = 1.024012134932392
$ws.Range("M3").Value = 1.03428910392326

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.020985824918356
$ws.Range("D4").Value = 1.025721490802771
$ws.Range("E4").Value = 1.021884833328113
$ws.Range("F4").Value = 1.032421622644335
$ws.Range("I4").Value = 1.029032540051098
$ws.Range("J4").Value = 1.025589802479126
$ws.Range("K4").Value = 1.028231729751108
$ws.Range("L4").Value = 1.024405025861917
$ws.Range("M4").Value = 1.03491467050195

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021242898067524
$ws.Range("D5").Value = 1.025902675849082
$ws.Range("E5").Value = 1.022103147084658
$ws.Range("F5").Value = 1.032737360771329
$ws.Range("I5").Value = 1.029071275866011
$ws.Range("J5").Value = 1.025747684513912
$ws.Range("K5").Value = 1.028360105661419
$ws.Range("L5").Value = 1.024570214453358
$ws.Range("M5").Value = 1.035177646557867

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.021286067247501
$ws.Range("D6").Value = 1.025933099272853
$ws.Range("E6").Value = 1.02213981300784
$ws.Range("F6").Value = 1.032790381909693
$ws.Range("I6").Value = 1.029077761456165
$ws.Range("J6").Value = 1.025774190863559
$ws.Range("K6").Value = 1.028381652691359
$ws.Range("L6").Value = 1.02459795133099
$ws.Range("M6").Value = 1.035221800677847

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.020989259575529
$ws.Range("D7").Value = 1.025723911694316
$ws.Range("E7").Value = 1.021887749770787
$ws.Range("F7").Value = 1.032425841062289
$ws.Range("I7").Value = 1.029033058868458
$ws.Range("J7").Value = 1.025591912289586
$ws.Range("K7").Value = 1.028233445641745
$ws.Range("L7").Value = 1.024407233053278
$ws.Range("M7").Value = 1.034918184448792

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.019749448109126
$ws.Range("D8").Value = 1.024849790277079
$ws.Range("E8").Value = 1.020835633889304
$ws.Range("F8").Value = 1.030903179197613
$ws.Range("I8").Value = 1.028843560451348
$ws.Range("J8").Value = 1.024829626520697
$ws.Range("K8").Value = 1.027612834058649
$ws.Range("L8").Value = 1.023610210471684
$ws.Range("M8").Value = 1.033649015758505

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.017567789615111
$ws.Range("D9").Value = 1.023310495848177
$ws.Range("E9").Value = 1.018987278934484
$ws.Range("F9").Value = 1.028224050016765
$ws.Range("I9").Value = 1.028499683461504
$ws.Range("J9").Value = 1.023484957438122
$ws.Range("K9").Value = 1.026515031778077
$ws.Range("L9").Value = 1.022206371813245
$ws.Range("M9").Value = 1.031412204358197

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01611544009761
$ws.Range("D10").Value = 1.022285058361629
$ws.Range("E10").Value = 1.017758880743581
$ws.Range("F10").Value = 1.026440649275957
$ws.Range("I10").Value = 1.02826376325615
$ws.Range("J10").Value = 1.022587598009859
$ws.Range("K10").Value = 1.02578038529293
$ws.Range("L10").Value = 1.021270944245853
$ws.Range("M10").Value = 1.029920753611832

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.015487060264339
$ws.Range("D11").Value = 1.021841227850909
$ws.Range("E11").Value = 1.017227895276424
$ws.Range("F11").Value = 1.025669052783165
$ws.Range("I11").Value = 1.028160030619841
$ws.Range("J11").Value = 1.022198823809673
$ws.Range("K11").Value = 1.025461626382033
$ws.Range("L11").Value = 1.02086601430635
$ws.Range("M11").Value = 1.029274882203542

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.015253727268612
$ws.Range("D12").Value = 1.02167639928823
$ws.Range("E12").Value = 1.017030802814988
$ws.Range("F12").Value = 1.025382541590321
$ws.Range("I12").Value = 1.028121263062009
$ws.Range("J12").Value = 1.022054384643775
$ws.Range("K12").Value = 1.025343127929341
$ws.Range("L12").Value = 1.020715623571476
$ws.Range("M12").Value = 1.029034967245629

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.015303774573711
$ws.Range("D13").Value = 1.021711754204382
$ws.Range("E13").Value = 1.017073073500944
$ws.Range("F13").Value = 1.025443994959025
$ws.Range("I13").Value = 1.028129589536972
$ws.Range("J13").Value = 1.022085368732402
$ws.Range("K13").Value = 1.025368550635779
$ws.Range("L13").Value = 1.020747882054983
$ws.Range("M13").Value = 1.029086430251622

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.015467771347671
$ws.Range("D14").Value = 1.021827602454616
$ws.Range("E14").Value = 1.017211600706393
$ws.Range("F14").Value = 1.025645367760226
$ws.Range("I14").Value = 1.028156830906541
$ws.Range("J14").Value = 1.022186885054022
$ws.Range("K14").Value = 1.025451833241344
$ws.Range("L14").Value = 1.02085358258393
$ws.Range("M14").Value = 1.029255050953476

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.015568825161567
$ws.Range("D15").Value = 1.02189898437645
$ws.Range("E15").Value = 1.017296970365817
$ws.Range("F15").Value = 1.025769452651867
$ws.Range("I15").Value = 1.028173583868055
$ws.Range("J15").Value = 1.022249428507163
$ws.Range("K15").Value = 1.025503133546488
$ws.Range("L15").Value = 1.020918710613529
$ws.Range("M15").Value = 1.029358942465373

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016157154119394
$ws.Range("D16").Value = 1.022314518051136
$ws.Range("E16").Value = 1.017794139958115
$ws.Range("F16").Value = 1.026491870802102
$ws.Range("I16").Value = 1.02827061444341
$ws.Range("J16").Value = 1.022613395271999
$ws.Range("K16").Value = 1.025801526598516
$ws.Range("L16").Value = 1.021297820631635
$ws.Range("M16").Value = 1.029963616631342

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.016526330818536
$ws.Range("D17").Value = 1.022575223248984
$ws.Range("E17").Value = 1.018106248146478
$ws.Range("F17").Value = 1.026945192544524
$ws.Range("I17").Value = 1.028331057053579
$ws.Range("J17").Value = 1.022841645823865
$ws.Range("K17").Value = 1.025988526629325
$ws.Range("L17").Value = 1.021535658044132
$ws.Range("M17").Value = 1.030342895762252

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.016741713264884
$ws.Range("D18").Value = 1.02272730649828
$ws.Range("E18").Value = 1.01828838417379
$ws.Range("F18").Value = 1.027209668035006
$ws.Range("I18").Value = 1.028366159961275
$ws.Range("J18").Value = 1.022974760042716
$ws.Range("K18").Value = 1.026097537609987
$ws.Range("L18").Value = 1.021674395757968
$ws.Range("M18").Value = 1.030564116748051

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.016815161224492
$ws.Range("D19").Value = 1.022779166026577
$ws.Range("E19").Value = 1.018350502846983
$ws.Range("F19").Value = 1.027299857626077
$ws.Range("I19").Value = 1.028378103310754
$ws.Range("J19").Value = 1.023020145059238
$ws.Range("K19").Value = 1.02613469683229
$ws.Range("L19").Value = 1.021721703604489
$ws.Range("M19").Value = 1.030639546351381

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.016486716691414
$ws.Range("D20").Value = 1.022547250145098
$ws.Range("E20").Value = 1.018072752719142
$ws.Range("F20").Value = 1.026896549146337
$ws.Range("I20").Value = 1.028324587878081
$ws.Range("J20").Value = 1.022817158811775
$ws.Range("K20").Value = 1.025968469803784
$ws.Range("L20").Value = 1.021510139179011
$ws.Range("M20").Value = 1.030302203333547

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.015419476303302
$ws.Range("D21").Value = 1.021793487198706
$ws.Range("E21").Value = 1.017170804046081
$ws.Range("F21").Value = 1.025586065893478
$ws.Range("I21").Value = 1.028148815531463
$ws.Range("J21").Value = 1.022156991881078
$ws.Range("K21").Value = 1.025427311266755
$ws.Range("L21").Value = 1.020822455901867
$ws.Range("M21").Value = 1.029205396623692

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014748895588266
$ws.Range("D22").Value = 1.021319739816483
$ws.Range("E22").Value = 1.016604518934477
$ws.Range("F22").Value = 1.024762657387572
$ws.Range("I22").Value = 1.028036931556466
$ws.Range("J22").Value = 1.021741739028542
$ws.Range("K22").Value = 1.025086501427516
$ws.Range("L22").Value = 1.020390188251337
$ws.Range("M22").Value = 1.028515735390594

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01510434162464
$ws.Range("D23").Value = 1.021570865440125
$ws.Range("E23").Value = 1.016904640635638
$ws.Range("F23").Value = 1.025199110292002
$ws.Range("I23").Value = 1.028096372998592
$ws.Range("J23").Value = 1.021961889229143
$ws.Range("K23").Value = 1.02526722418294
$ws.Range("L23").Value = 1.020619331188764
$ws.Range("M23").Value = 1.028881343113813

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.016504616462199
$ws.Range("D24").Value = 1.022559889930168
$ws.Range("E24").Value = 1.018087887587293
$ws.Range("F24").Value = 1.026918528816746
$ws.Range("I24").Value = 1.028327511490457
$ws.Range("J24").Value = 1.022828223501777
$ws.Range("K24").Value = 1.025977532814541
$ws.Range("L24").Value = 1.021521670021237
$ws.Range("M24").Value = 1.030320590509234

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.018131435311792
$ws.Range("D25").Value = 1.023708311985212
$ws.Range("E25").Value = 1.019464451429344
$ws.Range("F25").Value = 1.0289161957523
$ws.Range("I25").Value = 1.02858976120023
$ws.Range("J25").Value = 1.023832750847977
$ws.Range("K25").Value = 1.026799333271335
$ws.Range("L25").Value = 1.022569219482204
$ws.Range("M25").Value = 1.031990517198356
